# chore: update Sheets via scheduled runner
# Refreshes cached market-price figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) for a handful of leve rows across several crafting-job
# sheets, matching the latest data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value2 = 5363.636
$ws.Range("I100").Value2 = 2833.3333
$ws.Range("J100").Value2 = 8400
$ws.Range("K100").Value2 = 2833.3333
$ws.Range("L100").Value2 = 8400
$ws.Range("M100").Value2 = -2292.3333
$ws.Range("N100").Value2 = -9482
$ws.Range("H113").Value2 = 2967.2222
$ws.Range("I113").Value2 = 2968.3333
$ws.Range("J113").Value2 = 2966.6667
$ws.Range("K113").Value2 = 2968.3333
$ws.Range("L113").Value2 = 2966.6667
$ws.Range("M113").Value2 = 285.6667000000002
$ws.Range("N113").Value2 = -9474.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 5466525
$ws.Range("I61").Value2 = 7753413.5
$ws.Range("J61").Value2 = 3402
$ws.Range("K61").Value2 = 7753413.5
$ws.Range("L61").Value2 = 3402
$ws.Range("M61").Value2 = -7753201.5
$ws.Range("N61").Value2 = -3826
$ws.Range("H74").Value2 = 2198.311
$ws.Range("I74").Value2 = 1509.0344
$ws.Range("J74").Value2 = 3447.625
$ws.Range("K74").Value2 = 1509.0344
$ws.Range("L74").Value2 = 3447.625
$ws.Range("M74").Value2 = -635.0344
$ws.Range("N74").Value2 = -5195.625
$ws.Range("H77").Value2 = 2198.311
$ws.Range("I77").Value2 = 1509.0344
$ws.Range("J77").Value2 = 3447.625
$ws.Range("K77").Value2 = 7545.172
$ws.Range("L77").Value2 = 17238.125
$ws.Range("M77").Value2 = -3177.172
$ws.Range("N77").Value2 = -25974.125
$ws.Range("H102").Value2 = 2407.5
$ws.Range("I102").Value2 = 2252
$ws.Range("K102").Value2 = 2252
$ws.Range("M102").Value2 = -630
$ws.Range("H110").Value2 = 900
$ws.Range("I110").Value2 = 825
$ws.Range("K110").Value2 = 825
$ws.Range("M110").Value2 = 1220
$ws.Range("H136").Value2 = 5466525
$ws.Range("I136").Value2 = 7753413.5
$ws.Range("J136").Value2 = 3402
$ws.Range("K136").Value2 = 23260240.5
$ws.Range("L136").Value2 = 10206
$ws.Range("M136").Value2 = -23257690.5
$ws.Range("N136").Value2 = -15306

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value2 = 390.75
$ws.Range("I94").Value2 = 362.875
$ws.Range("J94").Value2 = 502.25
$ws.Range("K94").Value2 = 362.875
$ws.Range("L94").Value2 = 502.25
$ws.Range("M94").Value2 = 88.125
$ws.Range("N94").Value2 = -1404.25
$ws.Range("H99").Value2 = 1413.75
$ws.Range("I99").Value2 = 1225
$ws.Range("J99").Value2 = 1602.5
$ws.Range("K99").Value2 = 1225
$ws.Range("L99").Value2 = 1602.5
$ws.Range("M99").Value2 = 273
$ws.Range("N99").Value2 = -4598.5
$ws.Range("H103").Value2 = 41593.6
$ws.Range("J103").Value2 = 41593.6
$ws.Range("L103").Value2 = 41593.6
$ws.Range("N103").Value2 = -43937.6
$ws.Range("H105").Value2 = 2397.5
$ws.Range("I105").Value2 = 2496.6667
$ws.Range("K105").Value2 = 2496.6667
$ws.Range("M105").Value2 = -749.6667000000002
$ws.Range("H107").Value2 = 126235
$ws.Range("I107").Value2 = 333900
$ws.Range("J107").Value2 = 1636
$ws.Range("K107").Value2 = 333900
$ws.Range("L107").Value2 = 1636
$ws.Range("M107").Value2 = -331980
$ws.Range("N107").Value2 = -5476
$ws.Range("H134").Value2 = 3367.95
$ws.Range("I134").Value2 = 3438.3635
$ws.Range("J134").Value2 = 3281.889
$ws.Range("K134").Value2 = 10315.0905
$ws.Range("L134").Value2 = 9845.667000000001
$ws.Range("M134").Value2 = -7780.0905
$ws.Range("N134").Value2 = -14915.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1698.2222
$ws.Range("I16").Value2 = 1852.75
$ws.Range("J16").Value2 = 1574.6
$ws.Range("K16").Value2 = 1852.75
$ws.Range("L16").Value2 = 1574.6
$ws.Range("M16").Value2 = -1565.75
$ws.Range("N16").Value2 = -2148.6
$ws.Range("H113").Value2 = 1698.2222
$ws.Range("I113").Value2 = 1852.75
$ws.Range("J113").Value2 = 1574.6
$ws.Range("K113").Value2 = 1852.75
$ws.Range("L113").Value2 = 1574.6
$ws.Range("M113").Value2 = 317.25
$ws.Range("N113").Value2 = -5914.6
$ws.Range("H122").Value2 = 1511.4166
$ws.Range("I122").Value2 = 1076.3684
$ws.Range("J122").Value2 = 1997.6471
$ws.Range("K122").Value2 = 3229.1052
$ws.Range("L122").Value2 = 5992.9413
$ws.Range("M122").Value2 = -779.1052
$ws.Range("N122").Value2 = -10892.9413

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 794.96826
$ws.Range("I5").Value2 = 630.5349
$ws.Range("J5").Value2 = 1148.5
$ws.Range("K5").Value2 = 1891.6047
$ws.Range("L5").Value2 = 3445.5
$ws.Range("M5").Value2 = -1779.6047
$ws.Range("N5").Value2 = -3669.5
$ws.Range("H12").Value2 = 73.478264
$ws.Range("I12").Value2 = 27.666666
$ws.Range("J12").Value2 = 159.375
$ws.Range("K12").Value2 = 82.99999800000001
$ws.Range("L12").Value2 = 478.125
$ws.Range("M12").Value2 = 90.00000199999999
$ws.Range("N12").Value2 = -824.125
$ws.Range("H14").Value2 = 100000330
$ws.Range("I14").Value2 = 100000330
$ws.Range("K14").Value2 = 300000990
$ws.Range("M14").Value2 = -300000817
$ws.Range("H63").Value2 = 6271.4287
$ws.Range("J63").Value2 = 8000
$ws.Range("L63").Value2 = 24000
$ws.Range("N63").Value2 = -25498
$ws.Range("H64").Value2 = 1335529.2
$ws.Range("I64").Value2 = 1225
$ws.Range("J64").Value2 = 1589682.5
$ws.Range("K64").Value2 = 3675
$ws.Range("L64").Value2 = 4769047.5
$ws.Range("M64").Value2 = -3405
$ws.Range("N64").Value2 = -4769587.5
$ws.Range("H66").Value2 = 6271.4287
$ws.Range("J66").Value2 = 8000
$ws.Range("L66").Value2 = 72000
$ws.Range("N66").Value2 = -79488
$ws.Range("H67").Value2 = 1335529.2
$ws.Range("I67").Value2 = 1225
$ws.Range("J67").Value2 = 1589682.5
$ws.Range("K67").Value2 = 3675
$ws.Range("L67").Value2 = 4769047.5
$ws.Range("M67").Value2 = -2739
$ws.Range("N67").Value2 = -4770919.5
$ws.Range("H103").Value2 = 1860
$ws.Range("I103").Value2 = 450
$ws.Range("J103").Value2 = 3975
$ws.Range("K103").Value2 = 1350
$ws.Range("L103").Value2 = 11925
$ws.Range("M103").Value2 = -471
$ws.Range("N103").Value2 = -13683
$ws.Range("H107").Value2 = 2637.3333
$ws.Range("J107").Value2 = 3677.4
$ws.Range("L107").Value2 = 11032.2
$ws.Range("N107").Value2 = -14872.2
$ws.Range("H135").Value2 = 794.96826
$ws.Range("I135").Value2 = 630.5349
$ws.Range("J135").Value2 = 1148.5
$ws.Range("K135").Value2 = 5674.8141
$ws.Range("L135").Value2 = 10336.5
$ws.Range("M135").Value2 = -3139.8141
$ws.Range("N135").Value2 = -15406.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 2105.4443
$ws.Range("I40").Value2 = 2107.1428
$ws.Range("J40").Value2 = 2099.5
$ws.Range("K40").Value2 = 2107.1428
$ws.Range("L40").Value2 = 2099.5
$ws.Range("M40").Value2 = -1971.1428
$ws.Range("N40").Value2 = -2371.5
$ws.Range("H100").Value2 = 41400.305
$ws.Range("I100").Value2 = 44927.76
$ws.Range("J100").Value2 = 4362
$ws.Range("K100").Value2 = 44927.76
$ws.Range("L100").Value2 = 4362
$ws.Range("M100").Value2 = -44386.76
$ws.Range("N100").Value2 = -5444

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value2 = 3800
$ws.Range("I96").Value2 = 3110
$ws.Range("K96").Value2 = 3110
$ws.Range("M96").Value2 = -1737
